# regen sval data to filter save games
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.505614041169197, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 2.770546300948285)
    3 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    4 = @(0.06328177979961902, 0.3375848360084654, 16.98373111632243, 0.4998867070740569, 17.88448443920457)
    5 = @(0.7287194209349384, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 5.964442013611383)
    6 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 6.48142807727062, 14.40014219143469)
    7 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 6.048734245549538)
    8 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 5.488907176552729)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 7).Value = $vals[4]
}
